$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - "from" headers (C1:E1 adopt the same header style as B1)
$ws.Range("B1").Value = "Q_from_net1"
$ws.Range("C1").Value = "Q_from_CHP1"
$ws.Range("D1").Value = "Q_from_solar_th1"
$ws.Range("E1").Value = "Q_from_pvt1"

$ws.Range("B1").Copy()
$ws.Range("C1:E1").PasteSpecial(-4122)

# Row 2 - "demand" row
$ws.Range("A2").Value = "Q_to_demand1"
$ws.Range("B2").Value = "Q_net1_demand1"
$ws.Range("C2").Value = "Q_CHP1_demand1"
$ws.Range("D2").Value = "Q_solar_th1_demand1"
$ws.Range("E2").Value = "Q_pvt1_demand1"

# Row 3 - "net" row (B3 keeps its numeric objective value)
$ws.Range("A3").Value = "Q_to_net1"
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = "Q_CHP1_net1"
$ws.Range("D3").Value = "Q_solar_th1_net1"
$ws.Range("E3").Value = "Q_pvt1_net1"
